$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column (D2:D49) to text format before writing new numeric-looking
# values, so Excel keeps them as text (matching the original inline-string cells),
# then restore the default "Normal" style so no stray style index is left on the cells.
$ws.Range("D2:D49").NumberFormat = "@"

$ws.Range("D2").Value = "27.461.99"
$ws.Range("D3").Value = "1.837.35"
$ws.Range("D4").Value = "1.005"
$ws.Range("D5").Value = "332.57"
$ws.Range("D6").Value = "1.004"
$ws.Range("D7").Value = "0.4606"
$ws.Range("D8").Value = "0.3796"
$ws.Range("D9").Value = "46.50"
$ws.Range("D10").Value = "0.07873"
$ws.Range("D11").Value = "0.9712"
$ws.Range("D12").Value = "20.96"
$ws.Range("D13").Value = "1.835.30"
$ws.Range("D14").Value = "5.885"
$ws.Range("D15").Value = "7.007"
$ws.Range("D16").Value = "1.005"
$ws.Range("D17").Value = "87.76"
$ws.Range("D18").Value = "0.06650"
$ws.Range("D19").Value = "0.00001028"
$ws.Range("D20").Value = "16.93"
$ws.Range("D22").Value = "27.446.99"
$ws.Range("D23").Value = "5.329"
$ws.Range("D24").Value = "10.77"
$ws.Range("D25").Value = "2.299"
$ws.Range("D26").Value = "157.34"
$ws.Range("D27").Value = "19.30"
$ws.Range("D28").Value = "2.062"
$ws.Range("D29").Value = "5.307"
$ws.Range("D30").Value = "118.59"
$ws.Range("D31").Value = "0.9509"
$ws.Range("D32").Value = "0.09282"
$ws.Range("D33").Value = "3.571"
$ws.Range("D34").Value = "5.216"
$ws.Range("D35").Value = "1.318"
$ws.Range("D36").Value = "0.05925"
$ws.Range("D37").Value = "0.02182"
$ws.Range("D38").Value = "8.062"
$ws.Range("D39").Value = "1.159"
$ws.Range("D40").Value = "0.5792"
$ws.Range("D41").Value = "0.1834"
$ws.Range("D42").Value = "10.00"
$ws.Range("D43").Value = "1.239"
$ws.Range("D44").Value = "0.5476"
$ws.Range("D45").Value = "11.95"
$ws.Range("D46").Value = "1.861"
$ws.Range("D47").Value = "0.06655"
$ws.Range("D48").Value = "109.60"
$ws.Range("D49").Value = "1.039"

# Restore default styling on the Price column so no extra style index remains
$ws.Range("D2:D49").Style = "Normal"

# Update the Volume(1h) column (text cells, Excel keeps these as text naturally
# because of the leading/trailing spaces)
$ws.Range("E2").Value = "  -2.21%  "
$ws.Range("E3").Value = "  -2.81%  "
$ws.Range("E4").Value = "  -0.84%  "
$ws.Range("E5").Value = "  -1.17%  "
$ws.Range("E6").Value = "  -0.80%  "
$ws.Range("E7").Value = "  -3.39%  "
$ws.Range("E8").Value = "  -4.31%  "
$ws.Range("E9").Value = "  -1.70%  "
$ws.Range("E10").Value = "  -2.29%  "
$ws.Range("E11").Value = "  -5.08%  "
$ws.Range("E12").Value = "  -4.65%  "
$ws.Range("E13").Value = "  -3.86%  "
$ws.Range("E14").Value = "  -2.79%  "
$ws.Range("E15").Value = "  -3.25%  "
$ws.Range("E16").Value = "  -1.01%  "
$ws.Range("E17").Value = "  -1.03%  "
$ws.Range("E18").Value = "  -1.75%  "
$ws.Range("E19").Value = "  -2.38%  "
$ws.Range("E20").Value = "  -1.19%  "
$ws.Range("E21").Value = "  -0.99%  "
$ws.Range("E22").Value = "  -2.21%  "
$ws.Range("E23").Value = "  -3.80%  "
$ws.Range("E24").Value = "  -2.68%  "
$ws.Range("E25").Value = "  -2.32%  "
$ws.Range("E26").Value = "  -1.95%  "
$ws.Range("E27").Value = "  -3.73%  "
$ws.Range("E28").Value = "  -2.55%  "
$ws.Range("E29").Value = "  -4.24%  "
$ws.Range("E31").Value = "  -3.21%  "
$ws.Range("E32").Value = "  -3.50%  "
$ws.Range("E33").Value = "  -1.89%  "
$ws.Range("E34").Value = "  -2.99%  "
$ws.Range("E35").Value = "  -3.96%  "
$ws.Range("E36").Value = "  -2.69%  "
$ws.Range("E37").Value = "  -3.48%  "
$ws.Range("E38").Value = "  -2.24%  "
$ws.Range("E39").Value = "  -4.03%  "
$ws.Range("E40").Value = "  -3.34%  "
$ws.Range("E41").Value = "  -3.51%  "
$ws.Range("E42").Value = "  -3.71%  "
$ws.Range("E43").Value = "  -1.83%  "
$ws.Range("E44").Value = "  -3.61%  "
$ws.Range("E45").Value = "  -2.70%  "
$ws.Range("E46").Value = "  -4.11%  "
$ws.Range("E47").Value = "  -2.59%  "
$ws.Range("E48").Value = "  -2.44%  "
$ws.Range("E49").Value = "  -3.37%  "

# Row 50: BabyDogeCoin entry dropped; PaxDollar moves up into row 50 with updated data
$ws.Range("B50").Value = "PaxDollar"
$ws.Range("C50").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.004"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.90%  "

# Row 51: a new Aave entry replaces the old PaxDollar row
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "69.41"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.61%  "
